$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "atul "
$ws.Range("G6").Value = "aj "

$ws.Range("G7").Select()
